$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column A (the styled duplicate GENE-number column) - this
# shifts columns B:F left to A:E for every row.
$ws.Range("A1").EntireColumn.Delete()
